# Adds 15 new "harmonized name" columns (L..Z) to the header row (row 15)
# of the DDBJ BioSample 'Human' submission template, matching the style of
# the existing optional (yellow) header fields, plus field-definition
# comments for the subset of new columns that have one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values for columns L through Z on row 15.
$newHeaders = @{
    "L" = "cell_line";
    "M" = "cell_subtype";
    "N" = "cell_type";
    "O" = "culture_collection";
    "P" = "dev_stage";
    "Q" = "disease";
    "R" = "disease_stage";
    "S" = "ethnicity";
    "T" = "health_state";
    "U" = "karyotype";
    "V" = "phenotype";
    "W" = "population";
    "X" = "race";
    "Y" = "sample_type";
    "Z" = "treatment";
}

$columnOrder = @("L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

foreach ($col in $columnOrder) {
    $ws.Range("$col`15").Value = $newHeaders[$col]
}

# Match the existing "optional field" (yellow) header formatting used by
# C15/E15/F15, by copying their cell format onto the new header cells.
$ws.Range("C15").Copy()
$ws.Range("L15:Z15").PasteSpecial(-4122)  # xlPasteFormats

# Field-definition comments, only on the columns that have one.
$newComments = @{
    "L" = "Name of the cell line.";
    "N" = "Type of cell of the sample or from which the sample was obtained.";
    "O" = "Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier";
    "P" = "Developmental stage at the time of sampling.";
    "Q" = "list of diseases diagnosed; can include multiple diagnoses. the value of the field depends on host; for humans the terms should be chosen from DO (Disease Ontology), free text for non-human. For DO terms, please see http://gemina.svn.sourceforge.net/viewvc/gemina/trunk/Gemina/ontologies/gemina_symptom.obo?view=log";
    "R" = "Stage of disease at the time of sampling.";
    "S" = "ethnicity of the subject";
    "T" = "Health or disease status of sample at time of collection";
    "V" = "Phenotype of sampled organism. For Phenotypic quality Ontology (PATO) (v1.269) terms, please see http://bioportal.bioontology.org/visualize/44601";
    "W" = "for human: ; for plants: filial generation, number of progeny, genetic structure";
    "Y" = "Sample type, such as cell culture, mixed culture, tissue sample, whole organism, single cell, metagenomic assembly";
}

$commentOrder = @("L","N","O","P","Q","R","S","T","V","W","Y")

foreach ($col in $commentOrder) {
    $ws.Range("$col`15").AddComment($newComments[$col]) | Out-Null
}
